$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.126.69"
$ws.Range("E2").Value = "  -0.87%  "
$ws.Range("D3").Value = "1.667.53"
$ws.Range("E3").Value = "  -1.25%  "
$ws.Range("E4").Value = "  -0.82%  "
$ws.Range("D5").Value = "'210.23"
$ws.Range("E5").Value = "  -3.88%  "
$ws.Range("D6").Value = "'0.5199"
$ws.Range("E6").Value = "  -4.61%  "
$ws.Range("E7").Value = "  -0.81%  "
$ws.Range("D8").Value = "'0.2637"
$ws.Range("E8").Value = "  -3.87%  "
$ws.Range("D9").Value = "'0.06233"
$ws.Range("E9").Value = "  -3.40%  "
$ws.Range("D10").Value = "'21.17"
$ws.Range("E10").Value = "  -3.67%  "
$ws.Range("D11").Value = "'0.07493"
$ws.Range("E11").Value = "  -2.51%  "
$ws.Range("D12").Value = "1.685.94"
$ws.Range("E12").Value = "  -0.12%  "
$ws.Range("D13").Value = "'4.423"
$ws.Range("E13").Value = "  -2.36%  "
$ws.Range("D14").Value = "'0.5600"
$ws.Range("E14").Value = "  -3.85%  "
$ws.Range("D15").Value = "'66.08"
$ws.Range("E15").Value = "  +1.31%  "
$ws.Range("D16").Value = "'0.000007912"
$ws.Range("E16").Value = "  -5.55%  "
$ws.Range("D17").Value = "26.164.89"
$ws.Range("E17").Value = "  -0.89%  "
$ws.Range("E18").Value = "  -0.80%  "
$ws.Range("D19").Value = "'4.785"
$ws.Range("E19").Value = "  -3.04%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "'186.77"
$ws.Range("E20").Value = "  -2.42%  "
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").Value = "'10.36"
$ws.Range("E21").Value = "  -5.53%  "
$ws.Range("D22").Value = "'6.172"
$ws.Range("E22").Value = "  -1.30%  "
$ws.Range("E23").Value = "  -0.80%  "
$ws.Range("D24").Value = "'147.82"
$ws.Range("E24").Value = "  -1.14%  "
$ws.Range("E25").Value = "  -5.96%  "
$ws.Range("D26").Value = "'7.564"
$ws.Range("E26").Value = "  -4.03%  "
$ws.Range("D27").Value = "'16.06"
$ws.Range("E27").Value = "  +2.09%  "
$ws.Range("D28").Value = "'0.06277"
$ws.Range("E28").Value = "  -1.12%  "
$ws.Range("D29").Value = "'1.362"
$ws.Range("E29").Value = "  -2.87%  "
$ws.Range("D30").Value = "'1.274"
$ws.Range("E30").Value = "  -4.04%  "
$ws.Range("D31").Value = "'3.473"
$ws.Range("E31").Value = "  -3.08%  "
$ws.Range("D32").Value = "'3.425"
$ws.Range("E32").Value = "  -4.73%  "
$ws.Range("D33").Value = "'1.624"
$ws.Range("E33").Value = "  -3.42%  "
$ws.Range("D34").Value = "'0.9965"
$ws.Range("E34").Value = "  -4.41%  "
$ws.Range("D35").Value = "'0.6035"
$ws.Range("E35").Value = "  -1.82%  "
$ws.Range("D36").Value = "'2.403"
$ws.Range("E36").Value = "  -0.37%  "
$ws.Range("D37").Value = "'2.705"
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("D38").Value = "'6.128"
$ws.Range("E38").Value = "  -2.30%  "
$ws.Range("D39").Value = "'0.01606"
$ws.Range("E39").Value = "  -1.56%  "
$ws.Range("D40").Value = "1.075.20"
$ws.Range("E40").Value = "  -3.87%  "
$ws.Range("D41").Value = "'0.8632"
$ws.Range("E41").Value = "  -1.69%  "
$ws.Range("D42").Value = "'1.004"
$ws.Range("E42").Value = "  -1.13%  "
$ws.Range("D43").Value = "'99.18"
$ws.Range("E43").Value = "  -2.56%  "
$ws.Range("E44").Value = "  -1.39%  "
$ws.Range("D45").Value = "'0.00000000109"
$ws.Range("E45").Value = "  +0.68%  "
$ws.Range("D46").Value = "'56.02"
$ws.Range("E46").Value = "  -2.52%  "
$ws.Range("D47").Value = "'1.005"
$ws.Range("E47").Value = "  -1.09%  "
$ws.Range("D49").Value = "'7.951"
$ws.Range("E49").Value = "  -2.94%  "
$ws.Range("D50").Value = "'0.4247"
$ws.Range("E50").Value = "  -1.33%  "
$ws.Range("D51").Value = "'5.939"
$ws.Range("E51").Value = "  -3.14%  "
